$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P mirrors column O's formatting for rows 4-14 (data table extended
# with a 2022 column). Copy format first, then set values, so the copied
# cellXfs/style indices match column O exactly instead of Excel fabricating
# brand-new style entries.

$srcRange = $ws.Range("O4:O14")
$dstRange = $ws.Range("P4:P14")
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)

$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 1
$ws.Range("P6").Value = "-"
$ws.Range("P7").Value = "-"
$ws.Range("P8").Value = "-"
$ws.Range("P9").Value = "-"
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = "-"
$ws.Range("P12").Value = 1
$ws.Range("P13").Value = "-"
$ws.Range("P14").Value = "-"

$ws.Range("O21:O22").Select()
